{"js": "// Apply the tracked-change style corrections to the \"texte a JM de classe\" paragraph.\n\n// 1) \"figure ,\" -> \"figure  ,\" (an extra space is inserted before the comma)\nlet results = context.document.body.search(\"figure ,\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\"figure  ,\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Italicize \"SimDisplay\"\nresults = context.document.body.search(\"SimDisplay\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].font.italic = true;\nawait context.sync();\n\n// 3) \"N\u0153uds(\" -> \"n\u0153uds (\" (lower-case the leading letter, add a space before the parenthesis)\nresults = context.document.body.search(\"il est possible de cr\u00e9er des N\u0153uds(\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\"il est possible de cr\u00e9er des n\u0153uds (\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Italicize \"Node\"\nresults = context.document.body.search(\"Node\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].font.italic = true;\nawait context.sync();\n\n// 5) Insert \"des ar\u00eates (Line)\" right after \"Node),\" and add a space before \"(Directions)\"\nresults = context.document.body.search(\"), des besoins en transport(Directions)\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\"), des ar\u00eates (Line) des besoins en transport (Directions)\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 6) \"Circuit\" -> \"Route\" (the method-holder class was renamed in the prose)\nresults = context.document.body.search(\"impl\u00e9ment\u00e9s dans Circuit afin\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\"impl\u00e9ment\u00e9s dans Route afin\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 7) Italicize \"isStation\"\nresults = context.document.body.search(\"isStation\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].font.italic = true;\nawait context.sync();\n\n// 8) Add \" dans son parcours\" after \"...dernier n\u0153ud parcouru\"\nresults = context.document.body.search(\"retourner l\u2019index du dernier n\u0153ud parcouru. La classe\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\"retourner l\u2019index du dernier n\u0153ud parcouru dans son parcours. La classe\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 9) Italicize \"SimTimer\"\nresults = context.document.body.search(\"SimTimer\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].font.italic = true;\nawait context.sync();\n\n// 10) Rewrite the SimTimer description\nresults = context.document.body.search(\"utilis\u00e9e par Simulation,\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\"dicte le rythme de la simulation et\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"figure ,\" -> \"figure  ,\" (an extra space is inserted before the comma)\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"figure ,\")) {\n    $rng.Text = \"figure  ,\"\n}\n\n# 2) Italicize \"SimDisplay\"\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"SimDisplay\")) {\n    $rng.Font.Italic = $true\n}\n\n# 3) \"N\u0153uds(\" -> \"n\u0153uds (\" (lower-case the leading letter, add a space before the parenthesis)\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"il est possible de cr\u00e9er des N\u0153uds(\")) {\n    $rng.Text = \"il est possible de cr\u00e9er des n\u0153uds (\"\n}\n\n# 4) Italicize \"Node\"\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"Node\")) {\n    $rng.Font.Italic = $true\n}\n\n# 5) Insert \"des ar\u00eates (Line)\" right after \"Node),\" and add a space before \"(Directions)\"\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"), des besoins en transport(Directions)\")) {\n    $rng.Text = \"), des ar\u00eates (Line) des besoins en transport (Directions)\"\n}\n\n# 6) \"Circuit\" -> \"Route\" (the method-holder class was renamed in the prose)\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"impl\u00e9ment\u00e9s dans Circuit afin\")) {\n    $rng.Text = \"impl\u00e9ment\u00e9s dans Route afin\"\n}\n\n# 7) Italicize \"isStation\"\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"isStation\")) {\n    $rng.Font.Italic = $true\n}\n\n# 8) Add \" dans son parcours\" after \"...dernier n\u0153ud parcouru\"\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"retourner l\u2019index du dernier n\u0153ud parcouru. La classe\")) {\n    $rng.Text = \"retourner l\u2019index du dernier n\u0153ud parcouru dans son parcours. La classe\"\n}\n\n# 9) Italicize \"SimTimer\"\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"SimTimer\")) {\n    $rng.Font.Italic = $true\n}\n\n# 10) Rewrite the SimTimer description\n$rng = $d.Content\n$rng.Find.MatchCase = $true\nif ($rng.Find.Execute(\"utilis\u00e9e par Simulation,\")) {\n    $rng.Text = \"dicte le rythme de la simulation et\"\n}\n"}
